# Automatische test-sync: 2025-08-14 21:40:50
# Appends the new mail-log row (row 29) to the "Logs" sheet, extends the
# conditional-formatting ranges that covered the old last row, and bumps
# the "Intern verzoek / Actie voor medewerker" tally on the "Dashboard"
# sheet from 20 to 21.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row -----------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A29").Value = "Demo inplannen"
$logs.Range("B29").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C29").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D29").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E29").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F29").Value = "2025-08-14 21:39:57"
$logs.Range("G29").Value = "Nee"
$logs.Range("H29").Value = "Ja"
$logs.Range("I29").Value = "Nee"
$logs.Range("J29").Value = "Nee"

# --- Extend conditional formatting sqref from row 28 to row 29 --------
$colsToExtend = "D", "G", "H", "I", "J"
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "28")
    $newRange = $logs.Range($col + "2:" + $col + "29")
    $count = $oldRange.FormatConditions.Count()
    for ($i = 1; $i -le $count; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the count for this category -----------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 21
